$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.677.77"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "3.113.37"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.82"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.68"
$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("D8").Value = "3.114.49"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  -0.29%  "

$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("E12").Value = "  +2.88%  "

$ws.Range("D13").Value = "3.654.14"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.131"
$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.26"
$ws.Range("E15").Value = "  +2.60%  "

$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("D17").Value = "57.776.80"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "3.122.01"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.83"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.08"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "336.36"
$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  +2.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.68"
$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "0.0₃0922"
$ws.Range("E28").Value = "  +2.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.57"
$ws.Range("E29").Value = "  +3.73%  "

$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.91"
$ws.Range("E34").Value = "  +0.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.69"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.66"
$ws.Range("E36").Value = "  +3.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").Value = "  +2.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.07"
$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0663"
$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.55"
$ws.Range("E41").Value = "  +14.36%  "

$ws.Range("D42").Value = "3.162.55"
$ws.Range("E42").Value = "  +1.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.690"
$ws.Range("E43").Value = "  +5.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.91"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.85"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D47").Value = "2.297.11"
$ws.Range("E47").Value = "  +1.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0259"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.985"
$ws.Range("E49").Value = "  +6.96%  "

$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.01"
$ws.Range("E51").Value = "  +2.46%  "
